$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the Post-Test column (C) values, completing the CLASS stress data
$ws.Range("C2").Value = "A little stressful"
$ws.Range("C3").Value = "Moderately stressful"
$ws.Range("C4").Value = "A little stressful"
$ws.Range("C5").Value = "Not stressful"
$ws.Range("C6").Value = "A little stressful"
$ws.Range("C7").Value = "A little stressful"

# Update the active selection to reflect the new last-used cell
$ws.Range("C8").Select()
